$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 4281339
$ws.Range("I51").Value = 18524118
$ws.Range("J51").Value = 8505
$ws.Range("K51").Value = 18524118
$ws.Range("L51").Value = 8505
$ws.Range("M51").Value = -18523634
$ws.Range("N51").Value = -9473
# Row 108
$ws.Range("H108").Value = 38930.5
$ws.Range("J108").Value = 38930.5
$ws.Range("L108").Value = 38930.5
$ws.Range("N108").Value = -46610.5
# Row 109
$ws.Range("H109").Value = 37929.668
$ws.Range("J109").Value = 37929.668
$ws.Range("L109").Value = 37929.668
$ws.Range("N109").Value = -40703.668
# Row 120
$ws.Range("H120").Value = 48564
$ws.Range("J120").Value = 48564
$ws.Range("L120").Value = 48564
$ws.Range("N120").Value = -58240
# Row 126
$ws.Range("H126").Value = 46765.332
$ws.Range("J126").Value = 46765.332
$ws.Range("L126").Value = 46765.332
$ws.Range("N126").Value = -56645.332
# Row 128
$ws.Range("H128").Value = 46092
$ws.Range("J128").Value = 46092
$ws.Range("L128").Value = 46092
$ws.Range("N128").Value = -56052
# Row 130
$ws.Range("H130").Value = 47386
$ws.Range("J130").Value = 47386
$ws.Range("L130").Value = 47386
$ws.Range("N130").Value = -57426
# Row 132
$ws.Range("H132").Value = 26713.895
$ws.Range("I132").Value = 4418.852
$ws.Range("J132").Value = 81438.09
$ws.Range("K132").Value = 13256.556
$ws.Range("L132").Value = 244314.27
$ws.Range("M132").Value = -10726.556
$ws.Range("N132").Value = -249374.27

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 80
$ws.Range("H80").Value = 38316.668
$ws.Range("J80").Value = 38316.668
$ws.Range("L80").Value = 38316.668
$ws.Range("N80").Value = -40312.668
# Row 83
$ws.Range("H83").Value = 38316.668
$ws.Range("J83").Value = 38316.668
$ws.Range("L83").Value = 114950.004
$ws.Range("N83").Value = -124934.004
# Row 107
$ws.Range("H107").Value = 38678.332
$ws.Range("J107").Value = 38678.332
$ws.Range("L107").Value = 38678.332
$ws.Range("N107").Value = -46358.332
# Row 109
$ws.Range("H109").Value = 39977
$ws.Range("J109").Value = 39977
$ws.Range("L109").Value = 39977
$ws.Range("N109").Value = -42751
# Row 117
$ws.Range("H117").Value = 47997
$ws.Range("J117").Value = 47997
$ws.Range("L117").Value = 47997
$ws.Range("N117").Value = -57175
# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
# Row 123
$ws.Range("H123").Value = 38392
$ws.Range("J123").Value = 38392
$ws.Range("L123").Value = 38392
$ws.Range("N123").Value = -48192
# Row 128
$ws.Range("H128").Value = 50374.668
$ws.Range("J128").Value = 50374.668
$ws.Range("L128").Value = 50374.668
$ws.Range("N128").Value = -60334.668
# Row 131
$ws.Range("H131").Value = 44686
$ws.Range("J131").Value = 44686
$ws.Range("L131").Value = 44686
$ws.Range("N131").Value = -54766
# Row 132
$ws.Range("H132").Value = 2885.45
$ws.Range("I132").Value = 1654.1666
$ws.Range("J132").Value = 4732.375
$ws.Range("K132").Value = 4962.4998
$ws.Range("L132").Value = 14197.125
$ws.Range("M132").Value = -2432.4998
$ws.Range("N132").Value = -19257.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 124
$ws.Range("H124").Value = 50992
$ws.Range("J124").Value = 50992
$ws.Range("L124").Value = 50992
$ws.Range("N124").Value = -60812
# Row 126
$ws.Range("H126").Value = 50768
$ws.Range("J126").Value = 50768
$ws.Range("L126").Value = 50768
$ws.Range("N126").Value = -60648
# Row 130
$ws.Range("H130").Value = 49178.332
$ws.Range("J130").Value = 49178.332
$ws.Range("L130").Value = 49178.332
$ws.Range("N130").Value = -59218.332
# Row 133
$ws.Range("H133").Value = 49000
$ws.Range("J133").Value = 49000
$ws.Range("L133").Value = 49000
$ws.Range("N133").Value = -59120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 49883.5
$ws.Range("J20").Value = 49883.5
$ws.Range("L20").Value = 49883.5
$ws.Range("N20").Value = -50355.5
# Row 30
$ws.Range("H30").Value = 49883.5
$ws.Range("J30").Value = 49883.5
$ws.Range("L30").Value = 49883.5
$ws.Range("N30").Value = -50065.5
# Row 116
$ws.Range("H116").Value = 49822.332
$ws.Range("J116").Value = 49822.332
$ws.Range("L116").Value = 49822.332
$ws.Range("N116").Value = -59000.332
# Row 118
$ws.Range("H118").Value = 48742
$ws.Range("J118").Value = 48742
$ws.Range("L118").Value = 48742
$ws.Range("N118").Value = -52056
# Row 128
$ws.Range("H128").Value = 49883.5
$ws.Range("J128").Value = 49883.5
$ws.Range("L128").Value = 49883.5
$ws.Range("N128").Value = -59843.5
# Row 132
$ws.Range("H132").Value = 58567.76
$ws.Range("I132").Value = 2005.8235
$ws.Range("J132").Value = 178761.88
$ws.Range("K132").Value = 6017.470499999999
$ws.Range("L132").Value = 536285.64
$ws.Range("M132").Value = -3487.470499999999
$ws.Range("N132").Value = -541345.64

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 110
$ws.Range("H110").Value = 47688
$ws.Range("J110").Value = 47688
$ws.Range("L110").Value = 47688
$ws.Range("N110").Value = -55868
# Row 130
$ws.Range("H130").Value = 53984
$ws.Range("J130").Value = 53984
$ws.Range("L130").Value = 53984
$ws.Range("N130").Value = -64024
# Row 132
$ws.Range("H132").Value = 2517.8157
$ws.Range("I132").Value = 1764.8077
$ws.Range("J132").Value = 4149.3335
$ws.Range("K132").Value = 5294.4231
$ws.Range("L132").Value = 12448.0005
$ws.Range("M132").Value = -2764.4231
$ws.Range("N132").Value = -17508.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 48644.5
$ws.Range("J36").Value = 48644.5
$ws.Range("L36").Value = 48644.5
$ws.Range("N36").Value = -49768.5
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 127
$ws.Range("H127").Value = 49563.668
$ws.Range("J127").Value = 49563.668
$ws.Range("L127").Value = 49563.668
$ws.Range("N127").Value = -59483.668
# Row 132
$ws.Range("H132").Value = 3163.625
$ws.Range("I132").Value = 2370.92
$ws.Range("K132").Value = 7112.76
$ws.Range("M132").Value = -4582.76
# Row 133
$ws.Range("H133").Value = 39713.145
$ws.Range("J133").Value = 39713.145
$ws.Range("L133").Value = 39713.145
$ws.Range("N133").Value = -44773.145
# Row 137
$ws.Range("H137").Value = 40916.668
$ws.Range("J137").Value = 40916.668
$ws.Range("L137").Value = 40916.668
$ws.Range("N137").Value = -51116.668

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 44425
$ws.Range("J16").Value = 44425
$ws.Range("L16").Value = 44425
$ws.Range("N16").Value = -45009
# Row 119
$ws.Range("H119").Value = 46226.668
$ws.Range("J119").Value = 46226.668
$ws.Range("L119").Value = 46226.668
$ws.Range("N119").Value = -55902.668
